$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 307.18981030783209
$ws.Range("D2").Value = 258.89379074917952
$ws.Range("E2").Value = 364.49533719557132
$ws.Range("F2").Value = 90
$ws.Range("G2").Value = 321.30670273188684
$ws.Range("H2").Value = 230.2605568909878
$ws.Range("I2").Value = 24.271593835809995
$ws.Range("J2").Value = 269.69202395841899
$ws.Range("K2").Value = 135.11352815559499
$ws.Range("L2").Value = 411.36323703841498
$ws.Range("M2").Value = 273.07956277706433
$ws.Range("N2").Value = 369.53384268670936
$ws.Range("C3").Value = 372.68161026742007
$ws.Range("D3").Value = 319.89860542203439
$ws.Range("E3").Value = 434.17376717938771
$ws.Range("F3").Value = 90
$ws.Range("G3").Value = 378.66382257554795
$ws.Range("H3").Value = 247.44143834381498
$ws.Range("I3").Value = 26.082617755819228
$ws.Range("J3").Value = 338.26945309657754
$ws.Range("K3").Value = 183.425462756554
$ws.Range("L3").Value = 558.56727619768799
$ws.Range("M3").Value = 326.83821666737072
$ws.Range("N3").Value = 430.48942848372519
$ws.Range("C4").Value = 382.11853340437693
$ws.Range("D4").Value = 329.90001435210377
$ws.Range("E4").Value = 442.60250748358544
$ws.Range("F4").Value = 90
$ws.Range("G4").Value = 382.31663635731769
$ws.Range("H4").Value = 286.82702033937073
$ws.Range("I4").Value = 30.234222625061761
$ws.Range("J4").Value = 298.2746738667185
$ws.Range("K4").Value = 172.972731633774
$ws.Range("L4").Value = 501.32010073781998
$ws.Range("M4").Value = 322.24188000518899
$ws.Range("N4").Value = 442.39139270944639
$ws.Range("C5").Value = 411.09214673670857
$ws.Range("D5").Value = 348.08821798464055
$ws.Range("E5").Value = 485.49977958763463
$ws.Range("F5").Value = 90
$ws.Range("G5").Value = 435.22762068490727
$ws.Range("H5").Value = 371.2848152466151
$ws.Range("I5").Value = 39.13685589380384
$ws.Range("J5").Value = 339.73423531647154
$ws.Range("K5").Value = 196.59419738297501
$ws.Range("L5").Value = 556.15024332354903
$ws.Range("M5").Value = 357.46352165827216
$ws.Range("N5").Value = 512.99171971154237
$ws.Range("C6").Value = 504.60092544304371
$ws.Range("D6").Value = 435.15404772633633
$ws.Range("E6").Value = 585.13093303019264
$ws.Range("F6").Value = 90
$ws.Range("G6").Value = 508.32609093450469
$ws.Range("H6").Value = 355.25677837304471
$ws.Range("I6").Value = 37.447352462415665
$ws.Range("J6").Value = 433.87664252341904
$ws.Range("K6").Value = 256.33666210491498
$ws.Range("L6").Value = 680.55502750315304
$ws.Range("M6").Value = 433.91899923878071
$ws.Range("N6").Value = 582.73318263022861
$ws.Range("C7").Value = 589.2677085670739
$ws.Range("D7").Value = 507.55483293132846
$ws.Range("E7").Value = 684.13580135660072
$ws.Range("F7").Value = 90
$ws.Range("G7").Value = 617.7824122299628
$ws.Range("H7").Value = 433.4468756492106
$ws.Range("I7").Value = 45.689312391175918
$ws.Range("J7").Value = 538.92190508386398
$ws.Range("K7").Value = 287.79837666053203
$ws.Range("L7").Value = 823.75523332093906
$ws.Range("M7").Value = 526.99872171360767
$ws.Range("N7").Value = 708.56610274631794
